# Restored from revision of admin on 11/03/2020 07:58:35 AM.TEST Author: admin. Type: SAVE.
#
# The only substantive change in this revision is the "Integer min" value
# for rule R20 (cell C10 on the Rules sheet), which drops from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
